{"js": "// Change: \"The third question can be better visualized with the help of a facet\n// wrapped scatterplot as follows:\" -> \"... with the help of a scatterplot as follows:\"\n// (i.e. drop the words \"facet wrapped \" before \"scatterplot\").\n\nconst body = context.document.body;\n\n// Use a precise, uniquely-identifying search string (the search API does a\n// literal / wildcard-free match by default) so we only touch the intended run.\nconst searchResults = body.search(\"facet wrapped scatterplot as follows:\", {\n  matchCase: true,\n  matchWholeWord: false\n});\nsearchResults.load(\"items/text\");\nawait context.sync();\n\nif (searchResults.items.length === 0) {\n  throw new Error(\"Could not find target text 'facet wrapped scatterplot as follows:' in the document body.\");\n}\n\n// Replace \"facet wrapped scatterplot as follows:\" with \"scatterplot as follows:\"\n// by inserting the replacement text directly over the matched range.\nfor (let i = 0; i < searchResults.items.length; i++) {\n  searchResults.items[i].insertText(\"scatterplot as follows:\", Word.InsertLocation.replace);\n}\n\nawait context.sync();\n", "ps1": "# Change: \"The third question can be better visualized with the help of a facet\n# wrapped scatterplot as follows:\" -> \"... with the help of a scatterplot as follows:\"\n# (i.e. drop the words \"facet wrapped \" before \"scatterplot\").\n\n$d = $word.ActiveDocument\n\n$oldText = \"The third question can be better visualized with the help of a facet wrapped scatterplot as follows:\"\n$newText = \"The third question can be better visualized with the help of a scatterplot as follows:\"\n\nforeach ($p in $d.Paragraphs) {\n    if ($p.Range.Text -like \"*$oldText*\") {\n        $p.Range.Text = $newText\n    }\n}\n"}
